$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: Profitable = FALSE
$ws.Range("B2").Value = $false

# Add SellPrice and Price Change % to row 2
$ws.Range("E2").Value = 78.63
$ws.Range("F2").Value = -0.61931243680486481

# Holding changes from TRUE to FALSE
$ws.Range("G2").Value = $false

# New row 3 with updated Principle value
$ws.Range("C3").Value = 9938.07

# Re-fit column widths to the new content (bestFit-style resize)
# (column G / 7 is unchanged from its original best-fit width, so it is left alone)
$ws.Columns.Item(1).ColumnWidth = 13.417
$ws.Columns.Item(2).ColumnWidth = 7.251
$ws.Columns.Item(3).ColumnWidth = 6.917
$ws.Columns.Item(4).ColumnWidth = 6.417
$ws.Columns.Item(5).ColumnWidth = 6.084
$ws.Columns.Item(6).ColumnWidth = 11.584
